$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190564513206482
$ws.Range("B1").Value = 2.239068984985352
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.28122878074646
$ws.Range("E1").Value = 1.215537786483765
